$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (title reflects the new "through" date).
$ws.Name = "Through 2021-11-02"

# Row 4 (February) - 2018 column group (H/I/J arrest_made/no_arrest_made/arrest_rate; K/L/M is 2018 group)
$ws.Range("K4").Value = 5
$ws.Range("L4").Value = 46
$ws.Range("M4").Value = 0.098

# Row 13 (November running total) - update label and all year columns
$ws.Range("A13").Value = "November (through 11-02)"
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 5
$ws.Range("H13").Value = 1
$ws.Range("I13").Value = 6
$ws.Range("J13").Value = 0.1429
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 7
$ws.Range("M13").Value = 0.2222
$ws.Range("O13").Value = 3
$ws.Range("R13").Value = 13
$ws.Range("U13").Value = 13

# Row 14 (Total) - updated sums/rates reflecting the new November data
$ws.Range("C14").Value = 228
$ws.Range("D14").Value = 0.1231
$ws.Range("F14").Value = 439
$ws.Range("G14").Value = 0.1059
$ws.Range("H14").Value = 62
$ws.Range("I14").Value = 655
$ws.Range("J14").Value = 0.0865
$ws.Range("L14").Value = 556
$ws.Range("M14").Value = 0.109
$ws.Range("O14").Value = 437
$ws.Range("P14").Value = 0.099
$ws.Range("R14").Value = 1016
$ws.Range("S14").Value = 0.0505
$ws.Range("U14").Value = 1372
$ws.Range("V14").Value = 0.0583
